$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# These cells hold the value as text (shared string) in the workbook, so we
# must force a text number format before assigning the new value - otherwise
# a numeric-looking string like "10.77" would be auto-converted to a number.
$updates = @{
    "B11" = "10.77";
    "C11" = "3.12";
    "D11" = "13.89";
    "B12" = "17.39";
    "C12" = "40.35";
    "D12" = "57.75";
    "B14" = "75.64";
    "C14" = "21.92";
    "D14" = "97.56";
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
}
